$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B8").Value = "Dump20160214-1"
$ws.Range("C8").Value = "Added new tables (UNIT_CONVERT)"
